# OLX Monitor update — 2026-02-22 11:36
# The "PODSUMOWANIE" sheet holds a running log of scraped OLX listings in
# A:H (rows 7..N). This run re-scraped the last 8 listings that were already
# present at the bottom of the log (rows 107-114) and appended them again as
# a fresh snapshot (rows 115-122), with only the check timestamp (col A) and
# the view counter (col F) changed — everything else (profile, title, price,
# publish date, URL, slug, and all cell formatting) is identical to the
# previous snapshot of the same listings.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PODSUMOWANIE")

# Duplicate the previous snapshot block (rows 107-114) into the new block
# (rows 115-122). Range.Copy brings along values AND formatting, so the
# styles already used for A/C/D/E/F in that block (s="13"/"14"/"15") are
# reproduced exactly without having to rebuild them cell by cell.
$src = $ws.Range("A107:H114")
$dst = $ws.Range("A115:H122")
$src.Copy($dst)

# New check timestamp for every row in the fresh snapshot.
$timestamp = "2026-02-22 11:36:29"
$ws.Range("A115:A122").Value = $timestamp

# Updated view counters (column F) for each of the 8 listings, in the same
# order as the source block (107->115 ... 114->122).
$views = @(64, 34, 116, 135, 195, 33, 519, 79)
for ($i = 0; $i -lt $views.Length; $i++) {
    $row = 115 + $i
    $ws.Cells.Item($row, 6).Value = $views[$i]
}
